$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on D-column cells whose new value looks like a plain number,
# so Excel keeps them as text (preserving formatting such as trailing zeros and
# the multi-dot "thousand separator" notation used in this sheet) instead of
# silently converting them to numeric values.
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"

$ws.Range("D2").Value = "60.770.91"
$ws.Range("E2").Value = "  -2.89%  "

$ws.Range("D3").Value = "2.905.18"
$ws.Range("E3").Value = "  -3.83%  "

$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.05%  "

$ws.Range("D5").Value = "587.65"
$ws.Range("E5").Value = "  -1.51%  "

$ws.Range("D6").Value = "146.25"
$ws.Range("E6").Value = "  -3.02%  "

$ws.Range("E7").Value = "  +0.03%  "

$ws.Range("D8").Value = "0.503"
$ws.Range("E8").Value = "  -2.91%  "

$ws.Range("D9").Value = "2.903.95"
$ws.Range("E9").Value = "  -3.80%  "

$ws.Range("D10").Value = "6.73"
$ws.Range("E10").Value = "  +5.01%  "

$ws.Range("D11").Value = "0.144"
$ws.Range("E11").Value = "  -4.52%  "

$ws.Range("E12").Value = "  -2.73%  "

$ws.Range("E13").Value = "  -4.05%  "

$ws.Range("D14").Value = "33.65"
$ws.Range("E14").Value = "  -2.80%  "

$ws.Range("E15").Value = "  +0.42%  "

$ws.Range("D16").Value = "3.387.82"
$ws.Range("E16").Value = "  -3.78%  "

$ws.Range("D17").Value = "60.693.11"
$ws.Range("E17").Value = "  -2.93%  "

$ws.Range("E18").Value = "  -3.15%  "

$ws.Range("D19").Value = "2.904.88"
$ws.Range("E19").Value = "  -3.83%  "

$ws.Range("D20").Value = "425.71"
$ws.Range("E20").Value = "  -5.22%  "

$ws.Range("D21").Value = "13.57"
$ws.Range("E21").Value = "  -4.37%  "

$ws.Range("D22").Value = "0.670"
$ws.Range("E22").Value = "  -3.03%  "

$ws.Range("D23").Value = "7.09"
$ws.Range("E23").Value = "  -5.04%  "

$ws.Range("B24").Value = "RenderToken"
$ws.Range("C24").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D24").Value = "11.06"
$ws.Range("E24").Value = "  +1.34%  "

$ws.Range("B25").Value = "Litecoin"
$ws.Range("C25").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D25").Value = "80.20"
$ws.Range("E25").Value = "  -2.64%  "

$ws.Range("E26").Value = "  -1.43%  "

$ws.Range("D27").Value = "11.87"
$ws.Range("E27").Value = "  -1.88%  "

$ws.Range("D29").Value = "1.00"
$ws.Range("E29").Value = "  +0.02%  "

$ws.Range("D30").Value = "7.22"
$ws.Range("E30").Value = "  -1.93%  "

$ws.Range("D31").Value = "2.18"
$ws.Range("E31").Value = "  +1.50%  "

$ws.Range("E32").Value = "  -3.39%  "

$ws.Range("D33").Value = "26.44"
$ws.Range("E33").Value = "  -4.01%  "

$ws.Range("E34").Value = "  -2.88%  "

$ws.Range("D35").Value = "0.0₃0841"
$ws.Range("E35").Value = "  -1.83%  "

$ws.Range("E36").Value = "  -2.07%  "

$ws.Range("E37").Value = "  -3.78%  "

$ws.Range("D38").Value = "2.97"
$ws.Range("E38").Value = "  -3.69%  "

$ws.Range("D39").Value = "49.34"
$ws.Range("E39").Value = "  -1.72%  "

$ws.Range("D40").Value = "2.03"
$ws.Range("E40").Value = "  -2.71%  "

$ws.Range("E41").Value = "  +0.76%  "

$ws.Range("D42").Value = "8.68"
$ws.Range("E42").Value = "  -3.99%  "

$ws.Range("E43").Value = "  +1.03%  "

$ws.Range("D44").Value = "41.37"
$ws.Range("E44").Value = "  +1.80%  "

$ws.Range("E45").Value = "  -2.29%  "

$ws.Range("D46").Value = "371.85"
$ws.Range("E46").Value = "  -5.55%  "

$ws.Range("D47").Value = "2.661.47"
$ws.Range("E47").Value = "  -2.83%  "

$ws.Range("D48").Value = "133.38"
$ws.Range("E48").Value = "  +0.63%  "

$ws.Range("B49").Value = "USDe"
$ws.Range("C49").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D49").Value = "1.00"
$ws.Range("E49").Value = "  -0.02%  "

$ws.Range("B50").Value = "InjectiveProtocol"
$ws.Range("C50").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D50").Value = "25.36"
$ws.Range("E50").Value = "  +6.13%  "

$ws.Range("E51").Value = "  -1.37%  "
